# Update cryptocurrency price/volume data per latest GitHub Actions scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'63.209.69"
$ws.Range("E2").Value = "  +0.56%  "

$ws.Range("D3").Value = "'2.570.47"
$ws.Range("E3").Value = "  +1.07%  "

$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").Value = "'585.04"
$ws.Range("E5").Value = "  +3.14%  "

$ws.Range("D6").Value = "'147.41"
$ws.Range("E6").Value = "  +0.42%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("E8").Value = "  +2.95%  "

$ws.Range("E9").Value = "  +3.83%  "

$ws.Range("E10").Value = "  +0.61%  "

$ws.Range("E11").Value = "  +0.05%  "

$ws.Range("E12").Value = "  +1.46%  "

$ws.Range("D13").Value = "'27.41"
$ws.Range("E13").Value = "  +0.95%  "

$ws.Range("D14").Value = "'3.030.97"
$ws.Range("E14").Value = "  +1.04%  "

$ws.Range("D15").Value = "'63.170.61"
$ws.Range("E15").Value = "  +0.48%  "

$ws.Range("D16").Value = "'0.0000147"
$ws.Range("E16").Value = "  +3.96%  "

$ws.Range("D17").Value = "'2.612.72"
$ws.Range("E17").Value = "  +2.43%  "

$ws.Range("E18").Value = "  -0.74%  "

$ws.Range("D19").Value = "'342.89"
$ws.Range("E19").Value = "  +2.26%  "

$ws.Range("E20").Value = "  +3.03%  "

$ws.Range("D21").Value = "'6.88"

$ws.Range("D22").Value = "'1.00"
$ws.Range("E22").Value = "  -0.02%  "

$ws.Range("D23").Value = "'66.84"
$ws.Range("E23").Value = "  +3.12%  "

$ws.Range("D24").Value = "'2.697.21"
$ws.Range("E24").Value = "  +0.92%  "

$ws.Range("E25").Value = "  +1.17%  "

$ws.Range("E26").Value = "  +1.46%  "

$ws.Range("D27").Value = "'8.17"
$ws.Range("E27").Value = "  +12.47%  "

$ws.Range("E28").Value = "  +1.36%  "

$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").Value = "'1.49"
$ws.Range("E29").Value = "  -0.30%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("E30").Value = "  +0.15%  "

$ws.Range("D31").Value = "'1.99"
$ws.Range("E31").Value = "  +7.76%  "

$ws.Range("D33").Value = "'464.31"
$ws.Range("E33").Value = "  +13.44%  "

$ws.Range("E34").Value = "  +3.84%  "

$ws.Range("D35").Value = "'175.98"
$ws.Range("E35").Value = "  -0.48%  "

$ws.Range("E36").Value = "  +2.51%  "

$ws.Range("E37").Value = "  +1.59%  "

$ws.Range("D38").Value = "'4.54"
$ws.Range("E38").Value = "  +4.34%  "

$ws.Range("E39").Value = "  +0.05%  "

$ws.Range("E40").Value = "  -0.69%  "

$ws.Range("E41").Value = "  +0.00%  "

$ws.Range("D42").Value = "'151.46"
$ws.Range("E42").Value = "  -1.19%  "

$ws.Range("E43").Value = "  +2.24%  "

$ws.Range("D44").Value = "'21.03"
$ws.Range("E44").Value = "  +1.42%  "

$ws.Range("E45").Value = "  +5.94%  "

$ws.Range("D46").Value = "'0.613"
$ws.Range("E46").Value = "  +1.44%  "

$ws.Range("D47").Value = "'0.0979"
$ws.Range("E47").Value = "  +2.34%  "

$ws.Range("E48").Value = "  +1.69%  "

$ws.Range("E49").Value = "  -1.22%  "

$ws.Range("E50").Value = "  -0.05%  "

$ws.Range("D51").Value = "'0.164"
$ws.Range("E51").Value = "  +3.99%  "
